$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("D6").Value = 45050
$ws.Range("L6").Value = 'Especial'
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 13000
$ws.Range("Q6").Value = '$/caja 18 kilos empedrada'
$ws.Range("S6").Value = 722

# Row 7
$ws.Range("D7").Value = 45050
$ws.Range("L7").Value = 'Primera'
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 12000
$ws.Range("S7").Value = 667

# Row 8
$ws.Range("D8").Value = 45020
$ws.Range("M8").Value = 60
$ws.Range("Q8").Value = '$/caja 18 kilos granel'

# Row 9
$ws.Range("D9").Value = 45044
$ws.Range("L9").Value = 'Especial'
$ws.Range("M9").Value = 40
$ws.Range("N9").Value = 13000
$ws.Range("O9").Value = 13000
$ws.Range("P9").Value = 13000
$ws.Range("Q9").Value = '$/caja 18 kilos empedrada'
$ws.Range("S9").Value = 722

# Row 10
$ws.Range("D10").Value = 45044
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 40
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 12000
$ws.Range("Q10").Value = '$/caja 18 kilos empedrada'
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 667
$ws.Range("T10").Value = 18

# Row 11
$ws.Range("D11").Value = 45021
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 12000
$ws.Range("P11").Value = 12000
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 667
$ws.Range("T11").Value = 18

# Row 12
$ws.Range("D12").Value = 44699
$ws.Range("M12").Value = 60
$ws.Range("Q12").Value = '$/caja 15 kilos granel'
$ws.Range("R12").Value = 'Provincia de Curicó'
$ws.Range("S12").Value = 867
$ws.Range("T12").Value = 15

# Row 13
$ws.Range("D13").Value = 44699
$ws.Range("M13").Value = 120
$ws.Range("N13").Value = 11000
$ws.Range("P13").Value = 11500
$ws.Range("Q13").Value = '$/caja 15 kilos granel'
$ws.Range("R13").Value = 'Provincia de Curicó'
$ws.Range("S13").Value = 767
$ws.Range("T13").Value = 15

# Row 14
$ws.Range("D14").Value = 45049

# Row 15
$ws.Range("D15").Value = 45049
$ws.Range("M15").Value = 60

# New row 16
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C16").Value = 'Ñuble'
$ws.Range("D16").Value = 45040
$ws.Range("D16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = 'Fruta'
$ws.Range("G16").Value = 100104
$ws.Range("H16").Value = 'Frutos de pepita'
$ws.Range("I16").Value = 100104003
$ws.Range("J16").Value = 'Membrillo'
$ws.Range("K16").Value = 'Champion'
$ws.Range("L16").Value = 'Especial'
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = 13000
$ws.Range("O16").Value = 13000
$ws.Range("P16").Value = 13000
$ws.Range("Q16").Value = '$/caja 18 kilos empedrada'
$ws.Range("R16").Value = "Región de O'Higgins"
$ws.Range("S16").Value = 722
$ws.Range("T16").Value = 18

# New row 17
$ws.Range("A17").Value = 7
$ws.Range("B17").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C17").Value = 'Ñuble'
$ws.Range("D17").Value = 45040
$ws.Range("D17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = 'Fruta'
$ws.Range("G17").Value = 100104
$ws.Range("H17").Value = 'Frutos de pepita'
$ws.Range("I17").Value = 100104003
$ws.Range("J17").Value = 'Membrillo'
$ws.Range("K17").Value = 'Champion'
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 40
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 12000
$ws.Range("Q17").Value = '$/caja 18 kilos empedrada'
$ws.Range("R17").Value = "Región de O'Higgins"
$ws.Range("S17").Value = 667
$ws.Range("T17").Value = 18
